$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Helpers: rebuild the run sequence of a whole paragraph via InsertXML so we
# get exact control over run boundaries (WordprocessingML run-coalescing on
# plain Range.Text/InsertAfter edits would otherwise merge adjacent runs
# that share identical formatting).
# ---------------------------------------------------------------------------

function XmlEscape($s) {
    return $s -replace '&', '&amp;' -replace '<', '&lt;' -replace '>', '&gt;'
}

# Build a single <w:r> element. $extra is raw XML inserted right before the
# <w:t> (e.g. <w:lastRenderedPageBreak/>).
function Build-Run($text, $rprXml, $extra) {
    $needPreserve = ($text -ne $text.Trim()) -or ($text -eq "") -or ($text.StartsWith(" ")) -or ($text.EndsWith(" "))
    $spaceAttr = ""
    if ($needPreserve) { $spaceAttr = ' xml:space="preserve"' }
    $esc = XmlEscape $text
    return '<w:r><w:rPr>' + $rprXml + '</w:rPr>' + $extra + '<w:t' + $spaceAttr + '>' + $esc + '</w:t></w:r>'
}

# Replace the run content of an entire paragraph (leaves <w:pPr> untouched)
# with the supplied array of raw <w:r>...</w:r> XML fragments.
function Replace-ParagraphRuns($para, [string[]]$runsXml) {
    $r = $para.Range
    $start = $r.Start
    $end = $r.End - 1   # stop before the paragraph mark
    if ($end -gt $start) {
        $target = $d.Range($start, $end)
        $target.Text = ""
    }
    $body = [string]::Join("", $runsXml)
    $xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
        '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
        '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
        '<pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
        '<w:body><w:p>' + $body + '</w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    $ins = $d.Range($start, $start)
    $ins.InsertXML($xml)
}

# ---------------------------------------------------------------------------
# Paragraph 1 - title (single run -> single run, safe as a plain replace)
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("Unraveling the Enigmatic Universe", $true, $false, $false, $false, $false, `
    $true, 1, $false, "The Enchanting Realm of Biology: Unveiling the Secrets of Life", 2) | Out-Null

# ---------------------------------------------------------------------------
# Paragraph 2 - byline name: "Ethan Cross" -> "Dr" + "." + " Jane Carter"
# ---------------------------------------------------------------------------
$rpr36 = '<w:rFonts w:ascii="Aptos" w:hAnsi="Aptos"/><w:color w:val="000000"/><w:sz w:val="36"/>'
$p2runs = @(
    (Build-Run "Dr" $rpr36 ""),
    (Build-Run "." $rpr36 ""),
    (Build-Run " Jane Carter" $rpr36 "")
)
Replace-ParagraphRuns $d.Paragraphs(2) $p2runs

# ---------------------------------------------------------------------------
# Paragraph 3 - email: "ethan" + "." + "cross@intelligentesia" + "." + "com"
#   -> "janecarter12@eduworld" + "." + "net"
# ---------------------------------------------------------------------------
$rpr32 = '<w:rFonts w:ascii="Aptos" w:hAnsi="Aptos"/><w:color w:val="000000"/><w:sz w:val="32"/>'
$p3runs = @(
    (Build-Run "janecarter12@eduworld" $rpr32 ""),
    (Build-Run "." $rpr32 ""),
    (Build-Run "net" $rpr32 "")
)
Replace-ParagraphRuns $d.Paragraphs(3) $p3runs

# ---------------------------------------------------------------------------
# Paragraph 5 - main body text (keeps its two <w:br/><w:br/> pairs)
# ---------------------------------------------------------------------------
$rpr24 = '<w:rFonts w:ascii="Aptos" w:hAnsi="Aptos"/><w:color w:val="000000"/><w:sz w:val="24"/>'
$brBr = '<w:br/></w:r><w:r><w:rPr>' + $rpr24 + '</w:rPr><w:br/>'

$p5runs = @(
    (Build-Run "Biology, the study of life, embarks us on a captivating journey into the intricate world of living organisms" $rpr24 ""),
    (Build-Run "." $rpr24 ""),
    (Build-Run " From the awe-inspiring complexity of the human body to the microscopic wonders of a single cell, biology unravels the mysteries of life's myriad forms and functions" $rpr24 ""),
    (Build-Run "." $rpr24 ""),
    (Build-Run " This enthralling subject not only unveils the mechanisms underlying our existence but also provides a profound understanding of the delicate balance between organisms and their environment" $rpr24 ""),
    (Build-Run "." $rpr24 $brBr),
    (Build-Run "In this introductory chapter, we delve into the fundamental concepts of biology, setting the stage for a deeper exploration of the fascinating world that awaits us" $rpr24 ""),
    (Build-Run "." $rpr24 ""),
    (Build-Run " We begin by examining the characteristics shared by all living organisms, delving into the essential processes that sustain life and the remarkable diversity that exists among species" $rpr24 ""),
    (Build-Run "." $rpr24 ""),
    (Build-Run " From the intricate workings of cells to the intricate interactions between organisms, we gain a glimpse into the profound interconnectedness of life forms" $rpr24 ""),
    (Build-Run "." $rpr24 $brBr),
    (Build-Run "Moving forward, we unravel the mysteries of DNA, the blueprint that governs the traits and characteristics of all living things" $rpr24 ""),
    (Build-Run "." $rpr24 ""),
    (Build-Run " We explore the intricate workings of genes, the basic units of heredity, and witness the astonishing precision of the genetic code" $rpr24 ""),
    (Build-Run "." $rpr24 ""),
    (Build-Run " Through hands-on experiments and real-world case studies, we investigate the intricate mechanisms of inheritance, variation, and evolution, unraveling the remarkable story of life's progression over millions of years" $rpr24 ""),
    (Build-Run "." $rpr24 "")
)
Replace-ParagraphRuns $d.Paragraphs(5) $p5runs

# ---------------------------------------------------------------------------
# Paragraph 7 - summary body text
# ---------------------------------------------------------------------------
$rprSum = '<w:rFonts w:ascii="Aptos" w:hAnsi="Aptos"/><w:color w:val="000000"/>'
$lastBreak = '<w:lastRenderedPageBreak/>'
$p7runs = @(
    (Build-Run "This introductory chapter sets the stage for a fascinating journey into the world of biology, exploring the fundamental principles of life and the incredible diversity of organisms" $rprSum ""),
    (Build-Run "." $rprSum ""),
    (Build-Run " It lays the groundwork for understanding the intricate mechanisms underlying life's processes, the interconnections between living organisms and their environment, and the remarkable story of evolution" $rprSum ""),
    (Build-Run "." $rprSum ""),
    (Build-Run " From the microscopic world of cells to the vast tapestry of ecosystems, " $rprSum ""),
    (Build-Run "biology unveils the secrets of life's enigmatic dance, inspiring us to appreciate the profound beauty and complexity of the living world" $rprSum $lastBreak),
    (Build-Run "." $rprSum "")
)
Replace-ParagraphRuns $d.Paragraphs(7) $p7runs

# ---------------------------------------------------------------------------
# Add a new empty paragraph right after the summary paragraph (before sectPr)
# ---------------------------------------------------------------------------
$p7 = $d.Paragraphs(7)
$p7.Range.InsertParagraphAfter() | Out-Null

Write-Output "done"
